$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (the first data row), pushing all existing
# data rows down by one. This mirrors the daily "prepend latest price" update.
$ws.Rows.Item(2).Insert()

# New row 2 gets the newest date (one day after the former top date) with
# the same price values used throughout the rest of the sheet. Temporarily
# force the date cell to Text format so the "yyyy-mm-dd" string is not
# auto-converted into a date serial number, then restore the cell's format
# to match the rest of the date column (plain, unstyled) via a formats-only
# paste from a neighboring row.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-01-28"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
